# Apply "Horarios actualizados Linea 141 - 557" update.
# New scrape timestamp
$nuevaHora = "04:52:35"

$wb = $excel.ActiveWorkbook

# --- Sheet 1: LP1912 ---
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: $nuevaHora"

$sheet1Data = @(
    @("04:53", "11_ETCHEVERRY", 1),
    @("05:17", "17_ROMERO", 25),
    @("05:22", "23_HERNANDEZ", 30),
    @("05:44", "14_ABASTO", 52),
    @("05:47", "17_ROMERO", 55),
    @("06:01", "16_SANTA ANA", 69),
    @("06:09", "10_OLMOS", 77),
    @("06:15", "215A_EL PATO", 83),
    @("06:30", "23_HERNANDEZ", 98),
    @("06:34", "11_ETCHEVERRY", 102),
    @("06:39", "17X38_ROMERO", 107),
    @("06:41", "16_SANTA ANA", 109)
)

$r = 6
foreach ($row in $sheet1Data) {
    $ws1.Cells.Item($r, 1).Value = $nuevaHora
    $ws1.Cells.Item($r, 2).Value = $row[0]
    $ws1.Cells.Item($r, 3).Value = $row[1]
    $ws1.Cells.Item($r, 4).Value = $row[2]
    $ws1.Cells.Item($r, 5).Value = "LP1912"
    $r = $r + 1
}

# --- Sheet 2: LP1912-215 ---
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: $nuevaHora"
$ws2.Cells.Item(6, 1).Value = $nuevaHora
$ws2.Cells.Item(6, 4).Value = 83

# --- Sheet 3: 6203-6173 ---
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: $nuevaHora"
